$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.760.46"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.078.86"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "58.35"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.383.61"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.78"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.91"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.094.27"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.717.29"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.05"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.37"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.75"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("E27").Value = "  +5.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.47"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.39"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.26%  "
$ws.Range("E31").Value = "  +2.83%  "
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0630"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.66"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.31"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0973"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.79"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0215"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.452.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.37"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.17"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.40"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.269.46"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.27%  "
